$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.526.11"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "1.474.11"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "'0.9534"
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").Value = "'277.41"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "'0.3610"
$ws.Range("E7").Value = "  -1.58%  "
$ws.Range("D8").Value = "'0.3056"
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("D9").Value = "'39.40"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").Value = "'1.058"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("D11").Value = "'0.06640"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "'5.508"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "'18.12"
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("D15").Value = "'6.182"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").Value = "'0.9540"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "1.474.32"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").Value = "'0.05935"
$ws.Range("E19").Value = "  +5.67%  "
$ws.Range("D20").Value = "'69.06"
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("D21").Value = "'5.492"
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("D24").Value = "'2.254"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "20.549.20"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").Value = "'143.08"
$ws.Range("E26").Value = "  +5.73%  "
$ws.Range("D27").Value = "'2.130"
$ws.Range("D28").Value = "'17.18"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").Value = "1.636.11"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").Value = "'113.76"
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("D31").Value = "'3.942"
$ws.Range("E31").Value = "  +4.86%  "
$ws.Range("D32").Value = "'4.998"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("D33").Value = "'0.07988"
$ws.Range("E33").Value = "  +4.35%  "
$ws.Range("D34").Value = "'0.8051"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("D36").Value = "'1.217"
$ws.Range("E36").Value = "  +7.16%  "
$ws.Range("D37").Value = "'0.05837"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("D38").Value = "'4.710"
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("D39").Value = "'0.02044"
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("D40").Value = "'0.9547"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("D41").Value = "'10.35"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "'0.1874"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("D43").Value = "'7.419"
$ws.Range("E43").Value = "  +4.55%  "
$ws.Range("D44").Value = "'0.5281"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").Value = "'3.518"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "'12.21"
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("D47").Value = "'118.16"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("D48").Value = "'0.5189"
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").Value = "'1.809"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").Value = "'0.06476"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").Value = "'0.9819"
$ws.Range("E51").Value = "  -0.78%  "
